# Updates to RMI data
# Remove the stale "last updated" timestamp that had been written into
# cell C1 of the "About" sheet (it was a date value, e.g. 44307 = 2021-04-21).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
$ws.Range("C1").Clear()
